$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.322
$ws.Range("A3").Value = -21.625
$ws.Range("C3").Value = -12.774
$ws.Range("E6").Value = 16.671
$ws.Range("C12").Value = -11.288
$ws.Range("A14").Value = -21.601
$ws.Range("E19").Value = 16.327
$ws.Range("A21").Value = -20.423
$ws.Range("A23").Value = -20.585
$ws.Range("C24").Value = -12.321
$ws.Range("E24").Value = 17.043
$ws.Range("A25").Value = -20.688
$ws.Range("B25").Value = 7.15
$ws.Range("C25").Value = -13.043
$ws.Range("A26").Value = -21.667
$ws.Range("B27").Value = 5.308
$ws.Range("A29").Value = -21.219
$ws.Range("E30").Value = 16.342
$ws.Range("B31").Value = 5.837999999999999
$ws.Range("E31").Value = 16.255
$ws.Range("E33").Value = 17.186
$ws.Range("B39").Value = 7.581999999999999
$ws.Range("E42").Value = 16.556
$ws.Range("B48").Value = 5.24
$ws.Range("C50").Value = -13.133
$ws.Range("B51").Value = 5.458999999999999
$ws.Range("B52").Value = 5.397
$ws.Range("A53").Value = -21.651
$ws.Range("C53").Value = -12.345
$ws.Range("B55").Value = 4.510000000000001
$ws.Range("E55").Value = 16.461
$ws.Range("B56").Value = 4.907
$ws.Range("A57").Value = -21.354
$ws.Range("B57").Value = 5.926
$ws.Range("C57").Value = -13.257
$ws.Range("E58").Value = 16.542
$ws.Range("A59").Value = -22.5
$ws.Range("C61").Value = -13.025
$ws.Range("C63").Value = -11.913
$ws.Range("E65").Value = 17.273
$ws.Range("A69").Value = -21.462
$ws.Range("C70").Value = -12.04
$ws.Range("E70").Value = 17.357
$ws.Range("B73").Value = 6.919
$ws.Range("E75").Value = 16.532
$ws.Range("A79").Value = -21.211
$ws.Range("A83").Value = -22.015
$ws.Range("E83").Value = 16.59
$ws.Range("C86").Value = -13.406
$ws.Range("E86").Value = 16.279
$ws.Range("B89").Value = 5.67
$ws.Range("B90").Value = 5.907
$ws.Range("A91").Value = -21.533
$ws.Range("B92").Value = 6.026999999999999
$ws.Range("A93").Value = -21.379
$ws.Range("E96").Value = 16.349
$ws.Range("E97").Value = 16.958
$ws.Range("C98").Value = -12.285
$ws.Range("C100").Value = -13.247
$ws.Range("C102").Value = -13.361
